# Update the build timestamp embedded in the "version" strings across the
# workbook, from "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST".
# The overall release label ("mines - January 30") itself stays the same;
# only the "(built on ...)" timestamp changes.

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$rangeA2 = $wsAbout.Range("A2")
$rangeA2.Value = $rangeA2.Value().Replace($oldStamp, $newStamp)

$rangeA6 = $wsAbout.Range("A6")
$rangeA6.Value = $rangeA6.Value().Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 27; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S
    $cell.Value = $cell.Value().Replace($oldStamp, $newStamp)
}
